$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new, blank column before column N ("Late"); this shifts the
# existing N:P ("Late" / "heading" / "Outstanding") columns one position to
# the right, becoming O:Q, and leaves the new N column empty.
$mWidth = $ws.Columns("M:M").ColumnWidth
$ws.Columns("N:N").Insert()

# A freshly inserted column inherits the width of the column to its left.
$ws.Columns("N:N").ColumnWidth = $mWidth

# Make the "Repayment schedule" sheet the active tab (it was "Input"
# before), and leave the selection on it at K12 (previously H8).
$ws.Activate()
$ws.Range("K12").Select() | Out-Null
